$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data rows (2-530)
$ws.Range("C2:C530").Value = 45184

# 2. Existing last row (530) gains an explicit row height (15pt, custom height)
$ws.Rows.Item(530).RowHeight = 15

# 3. Append new row 531
$ws.Cells.Item(531,1).Value = "A 42562-2023"
$ws.Cells.Item(531,2).Value = 45181
$ws.Cells.Item(531,3).Value = 45184
$ws.Cells.Item(531,4).Value = "NORRBOTTENS LÄN"
$ws.Cells.Item(531,5).Value = "GÄLLIVARE"
$ws.Cells.Item(531,6).Value = "Sveaskog"
$ws.Cells.Item(531,7).Value = 32.4
$ws.Cells.Item(531,8).Value = 0
$ws.Cells.Item(531,9).Value = 0
$ws.Cells.Item(531,10).Value = 0
$ws.Cells.Item(531,11).Value = 0
$ws.Cells.Item(531,12).Value = 0
$ws.Cells.Item(531,13).Value = 0
$ws.Cells.Item(531,14).Value = 0
$ws.Cells.Item(531,15).Value = 0
$ws.Cells.Item(531,16).Value = 0
$ws.Cells.Item(531,17).Value = 0
$ws.Cells.Item(531,18).Value = ""
$ws.Range("B531:C531").NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(531,18).WrapText = $true
$ws.Rows.Item(531).RowHeight = 15

# 4. Append new row 532
$ws.Cells.Item(532,1).Value = "A 42568-2023"
$ws.Cells.Item(532,2).Value = 45181
$ws.Cells.Item(532,3).Value = 45184
$ws.Cells.Item(532,4).Value = "NORRBOTTENS LÄN"
$ws.Cells.Item(532,5).Value = "GÄLLIVARE"
$ws.Cells.Item(532,6).Value = "Sveaskog"
$ws.Cells.Item(532,7).Value = 24.6
$ws.Cells.Item(532,8).Value = 0
$ws.Cells.Item(532,9).Value = 0
$ws.Cells.Item(532,10).Value = 0
$ws.Cells.Item(532,11).Value = 0
$ws.Cells.Item(532,12).Value = 0
$ws.Cells.Item(532,13).Value = 0
$ws.Cells.Item(532,14).Value = 0
$ws.Cells.Item(532,15).Value = 0
$ws.Cells.Item(532,16).Value = 0
$ws.Cells.Item(532,17).Value = 0
$ws.Cells.Item(532,18).Value = ""
$ws.Range("B532:C532").NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(532,18).WrapText = $true
